$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("166:166").Insert()

$ws.Range("A166").Value = 5
$ws.Range("B166").Value = "Macroferia Regional de Talca"
$ws.Range("C166").Value = "Maule"
$ws.Range("D166").Value = 44767
$ws.Range("E166").Value = 7
$ws.Range("F166").Value = 100112008
$ws.Range("G166").Value = "Coliflor"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 3000
$ws.Range("K166").Value = 1000
$ws.Range("L166").Value = 1000
$ws.Range("M166").Value = 1000
$ws.Range("N166").Value = "$/unidad"
$ws.Range("O166").Value = "Región del Maule"
$ws.Range("P166").Value = 1000
$ws.Range("Q166").Value = 1
$ws.Range("R166").Value = "Hortaliza"
